# "sửa lại phần quay thưởng" - add a new "Ngày thực hiện" / {{Exports.Date}}
# column (F) to the LuckyNumber report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell F1 ------------------------------------------------
# Match the formatting already used by the other header cells (e.g. E1):
# bold red font on yellow fill with a thin border, centered.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F1").Value = "Ngày thực hiện"

# --- New data cell F2 ----------------------------------------------------
# Match the formatting already used by the other data cells (e.g. D2):
# bordered, centered both ways.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F2").Value = "{{Exports.Date}}"

$excel.CutCopyMode = 0

# --- Column width for the new column -------------------------------------
# Excel quantizes stored column widths to whole-pixel boundaries for the
# sheet's default font, so the raw target width of 19.140625 characters
# cannot be represented exactly via the ColumnWidth property; 18.33 is the
# input that lands on the closest achievable stored width.
$ws.Columns.Item(6).ColumnWidth = 18.33

# --- Update the saved selection/active cell -------------------------------
$ws.Range("E11").Select() | Out-Null
